$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.1523570014496575, 0.5952141443068003, 0.9456359972335984, 0.9724381714194473, 0.972070713281458, 42)
    3  = @(0.06653507977193379, 0.7435833224688242, 0.9402494797260309, 0.969664622292693, 0.9793968179231123, 41)
    4  = @(-0.176468453233442, 0.9424928345913857, 1.732816681911777, 1.316364950122791, 1.321101159551913, 40)
    5  = @(0.03757363485166248, 0.8767392780908844, 1.179621929259555, 1.086104014015027, 1.099643433808863, 39)
    6  = @(-0.1919880024426686, 0.9766461047393071, 1.673456267335123, 1.293621377117402, 1.296467979630036, 38)
    7  = @(0.0107767997713723, 0.8215764532688238, 1.166133474069657, 1.079876601315936, 1.094717637139123, 37)
    8  = @(-0.2391464184133537, 0.8937093494297598, 1.381707013101544, 1.175460340930967, 1.167201489911635, 36)
    9  = @(-0.07625812354922183, 0.8357403131286991, 1.353152740845836, 1.163250936318487, 1.177694786940479, 35)
    10 = @(-0.1780953643465996, 0.9948789155006484, 1.498614834658811, 1.224179249398882, 1.229369064332424, 34)
    11 = @(-0.1309683401722337, 0.9196742041917211, 1.399190577028793, 1.182873863532707, 1.193828528042026, 33)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
